$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Fred VanVleet",          "PG",       "Houston Rockets"),
    @("Dillon Brooks",          "SG,SF",    "Houston Rockets"),
    @("Bennedict Mathurin",     "SG,SF",    "Indiana Pacers"),
    @("Aaron Gordon",           "PF,C",     "Denver Nuggets"),
    @("Jayson Tatum",           "SF,PF",    "Boston Celtics"),
    @("Tyus Jones",             "PG",       "Phoenix Suns"),
    @("Ivica Zubac",            "C",        "LA Clippers"),
    @("Bobby Portis",           "PF,C",     "Milwaukee Bucks"),
    @("Anthony Edwards",        "SG,SF",    "Minnesota Timberwolves"),
    @("Jaren Jackson Jr.",      "PF,C",     "Memphis Grizzlies"),
    @("James Harden",           "PG,SG",    "LA Clippers"),
    @("Anfernee Simons",        "PG,SG",    "Portland Trail Blazers"),
    @("Paul George",            "SG,SF,PF", "Philadelphia 76ers"),
    @("Andrew Wiggins",         "SF,PF",    "Golden State Warriors"),
    @("Jaden Ivey",             "PG,SG",    "Detroit Pistons"),
    @("Zion Williamson",        "PF,C",     "New Orleans Pelicans"),
    @("Kyle Kuzma",             "PF",       "Washington Wizards"),
    @("Giannis Antetokounmpo",  "PF,C",     "Milwaukee Bucks")
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $rowIndex++
}
